$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 93.39526366666666
$ws.Range("H2").Value = 280.185791
$ws.Range("I2").Value = 0.2167755775732346
$ws.Range("J2").Value = 0.2167755775732346
$ws.Range("M2").Value = 72.07569866666667
$ws.Range("N2").Value = 216.227096
$ws.Range("O2").Value = 0.4479522040449755
$ws.Range("P2").Value = 0.4479522040449755
$ws.Range("Q2").Value = 6731.528880932548
$ws.Range("R2").Value = 60583.75992839294
$ws.Range("S2").Value = 0.097105097757053
$ws.Range("T2").Value = 0.09710509775705299

$ws.Range("G3").Value = 93.39526366666666
$ws.Range("H3").Value = 280.185791
$ws.Range("I3").Value = 0.2167755775732346
$ws.Range("J3").Value = 0.2167755775732346
$ws.Range("O3").Value = 0.04737448730867841
$ws.Range("P3").Value = 0.0473744873086784
$ws.Range("Q3").Value = 711.9124019439416
$ws.Range("R3").Value = 6407.211617495474
$ws.Range("S3").Value = 0.01026963184857463
$ws.Range("T3").Value = 0.01026963184857463

$ws.Range("G4").Value = 93.39526366666666
$ws.Range("H4").Value = 280.185791
$ws.Range("I4").Value = 0.2167755775732346
$ws.Range("J4").Value = 0.2167755775732346
$ws.Range("O4").Value = 0.5046733086463462
$ws.Range("P4").Value = 0.5046733086463461
$ws.Range("Q4").Value = 7583.896053891444
$ws.Range("R4").Value = 68255.064485023
$ws.Range("S4").Value = 0.109400847967607
$ws.Range("T4").Value = 0.1094008479676069

$ws.Range("I5").Value = 0.5566060939249745
$ws.Range("J5").Value = 0.5566060939249745
$ws.Range("M5").Value = 72.07569866666667
$ws.Range("N5").Value = 216.227096
$ws.Range("O5").Value = 0.4479522040449755
$ws.Range("P5").Value = 0.4479522040449755
$ws.Range("Q5").Value = 17284.28099928929
$ws.Range("R5").Value = 155558.5289936036
$ws.Range("S5").Value = 0.249332926558557
$ws.Range("T5").Value = 0.249332926558557

$ws.Range("I6").Value = 0.5566060939249745
$ws.Range("J6").Value = 0.5566060939249745
$ws.Range("O6").Value = 0.04737448730867841
$ws.Range("P6").Value = 0.0473744873086784
$ws.Range("S6").Value = 0.02636892833258177
$ws.Range("T6").Value = 0.02636892833258176

$ws.Range("I7").Value = 0.5566060939249745
$ws.Range("J7").Value = 0.5566060939249745
$ws.Range("O7").Value = 0.5046733086463462
$ws.Range("P7").Value = 0.5046733086463461
$ws.Range("S7").Value = 0.2809042390338358
$ws.Range("T7").Value = 0.2809042390338358

$ws.Range("I8").Value = 0.226618328501791
$ws.Range("J8").Value = 0.2266183285017909
$ws.Range("M8").Value = 72.07569866666667
$ws.Range("N8").Value = 216.227096
$ws.Range("O8").Value = 0.4479522040449755
$ws.Range("P8").Value = 0.4479522040449755
$ws.Range("Q8").Value = 7037.175683423567
$ws.Range("R8").Value = 63334.5811508121
$ws.Range("S8").Value = 0.1015141797293656
$ws.Range("T8").Value = 0.1015141797293656

$ws.Range("I9").Value = 0.226618328501791
$ws.Range("J9").Value = 0.2266183285017909
$ws.Range("O9").Value = 0.04737448730867841
$ws.Range("P9").Value = 0.0473744873086784
$ws.Range("S9").Value = 0.01073592712752201
$ws.Range("T9").Value = 0.01073592712752201

$ws.Range("I10").Value = 0.226618328501791
$ws.Range("J10").Value = 0.2266183285017909
$ws.Range("O10").Value = 0.5046733086463462
$ws.Range("P10").Value = 0.5046733086463461
$ws.Range("S10").Value = 0.1143682216449034
$ws.Range("T10").Value = 0.1143682216449034
